$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Generated: 2026-02-15" "Generated: 2026-02-22"

Replace-Text "Total Federal Climate Resilience Investment: `$1,865,595,124 across 48 Tribal Nations" `
             "Total Federal Climate Resilience Investment: `$2,137,232,918 across 51 Tribal Nations"

Replace-Text "Aggregate Economic Impact: `$3,633,663,879 to `$4,844,885,171" `
             "Aggregate Economic Impact: `$4,114,430,909 to `$5,485,907,878"

Replace-Text "Estimated Jobs Supported: 16,150 to 30,281" "Estimated Jobs Supported: 18,286 to 34,287"

Replace-Text "Total Federal Climate Resilience Awards: `$1,865,595,124" `
             "Total Federal Climate Resilience Awards: `$2,137,232,918"

Replace-Text "Tribal Nations with Awards: 48 of 55 (87%)" "Tribal Nations with Awards: 51 of 55 (93%)"

Replace-Text "Investment Gap: 7 Tribal Nation(s) in this region have received zero federal climate resilience funding through tracked programs." `
             "Investment Gap: 4 Tribal Nation(s) in this region have received zero federal climate resilience funding through tracked programs."
